$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date for every existing data row (2-400)
#    from 45192 (2023-09-23) to 45202 (2023-10-03).
$ws.Range("C2:C400").Value = 45202

# 2. Row 400 gains an explicit custom row height of 15 (as seen in the diff).
$ws.Rows.Item(400).RowHeight = 15

# 3. Append eight new rows (401-408) of data.
$newRows = @(
    @{ Row = 401; A = "A 45713-2023"; B = 45195; C = 45202; G = 4 },
    @{ Row = 402; A = "A 46063-2023"; B = 45196; C = 45202; G = 3.1 },
    @{ Row = 403; A = "A 46648-2023"; B = 45197; C = 45202; G = 1.5 },
    @{ Row = 404; A = "A 46458-2023"; B = 45197; C = 45202; G = 1.6 },
    @{ Row = 405; A = "A 47002-2023"; B = 45201; C = 45202; G = 1.3 },
    @{ Row = 406; A = "A 47100-2023"; B = 45201; C = 45202; G = 1.3 },
    @{ Row = 407; A = "A 47049-2023"; B = 45201; C = 45202; G = 24.6 },
    @{ Row = 408; A = "A 47103-2023"; B = 45201; C = 45202; G = 6 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("B$row`:C$row").NumberFormat = "YYYY-MM-DD"
    $ws.Range("D$row").Value = "VÄSTERBOTTENS LÄN"
    $ws.Range("E$row").Value = "MALÅ"
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row`:Q$row").Value = 0
    $ws.Range("R$row").Value = ""
    $ws.Range("R$row").WrapText = $true

    # Rows 401-407 get the explicit 15pt custom row height; row 408 does not.
    if ($row -lt 408) {
        $ws.Rows.Item($row).RowHeight = 15
    }
}
